# Apply the "resolved backlog" update to the ITI sheet of the Backlog workbook.
#
# Summary of the edit:
#   1. Every tracked incident (rows 2-45, column I "Status") moves from
#      "Pendente" to "Resolvido", and gets highlighted with a yellow fill
#      so the change stands out.
#   2. All of the rows that the old "Responsavel" filter (Felipe Nascimento)
#      had hidden are unhidden again - the whole backlog is visible now.
#   3. The AutoFilter's "Felipe Nascimento" criterion on the Responsavel
#      column is cleared, so the filter dropdown no longer restricts rows.
#   4. The sheet's last active selection moves to L48.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")

# --- 1. Update Status column (I2:I45): Pendente -> Resolvido, highlight yellow ---
$statusRange = $ws.Range("I2:I45")
$statusRange.Value = "Resolvido"
$statusRange.Interior.Color = 65535   # RGB(255,255,0) = yellow

# --- 2. Unhide every data row so the whole backlog shows again ---
for ($r = 2; $r -le 45; $r++) {
    $ws.Rows.Item($r).Hidden = $false
}

# --- 3. Remove the active filter criteria on the Responsavel column (col 2) ---
[void]$ws.Range("A1:I45").AutoFilter(2)

# --- 4. Update the sheet's last selected cell ---
[void]$ws.Range("L48").Select()
